{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\n// Find the paragraph that ends the \"feedback\" sentence \u2014 new content is\n// inserted right after it.\nconst anchorText = \"There will be a page for feedback about the site.\";\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === anchorText) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\n\nif (anchor) {\n  // Insert in reverse order, each time right \"After\" the anchor paragraph,\n  // so the final reading order is: anchor, twitter/facebook line, hosting line.\n  anchor.insertParagraph(\n    \"Web hosting offer pages should be added.\",\n    Word.InsertLocation.after\n  );\n  anchor.insertParagraph(\n    \"Twitter and facebook links should be added so that people can like, share and comment.\",\n    Word.InsertLocation.after\n  );\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the paragraph ending the \"feedback\" sentence \u2014 the two new\n# requirement lines are inserted directly after it.\n$anchorText = \"There will be a page for feedback about the site.\"\n$anchor = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -eq ($anchorText + \"`r\")) {\n        $anchor = $p\n        break\n    }\n}\n\nif ($anchor -ne $null) {\n    # First new paragraph, right after the anchor.\n    $anchor.Range.InsertParagraphAfter()\n    $p1 = $anchor.Next()\n    $p1.Range.Text = \"Twitter and facebook links should be added so that people can like, share and comment.\"\n\n    # Second new paragraph, right after the first.\n    $p1.Range.InsertParagraphAfter()\n    $p2 = $p1.Next()\n    $p2.Range.Text = \"Web hosting offer pages should be added.\"\n}\n"}
